# Junction_Flooding_405 style edit: widen several data columns, replace the
# four data rows with a newer sample window, and drop the old trailing row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Widen columns (offset -0.8333333333333334 compensates for the
#     "character width -> stored width" padding this host applies). ---
$MDW_OFFSET = 0.8333333333333334
$colWidths = @{
    "B" = 8
    "C" = 8
    "E" = 8
    "F" = 8
    "G" = 8
    "H" = 8
    "I" = 8
    "J" = 8
    "K" = 8
    "L" = 8
    "M" = 8
    "O" = 8
    "P" = 8
    "Q" = 8
    "S" = 7
    "T" = 9
    "U" = 8
    "V" = 8
    "W" = 8
    "X" = 8
    "Y" = 7
    "Z" = 8
    "AA" = 8
    "AB" = 8
    "AC" = 8
    "AD" = 8
    "AF" = 8
}
foreach ($col in $colWidths.Keys) {
    $ws.Range($col + ":" + $col).ColumnWidth = $colWidths[$col] - $MDW_OFFSET
}

# --- 2. Overwrite rows 2-5 with the new data sample. ---
$rowData = @{
    2 = @(45086.50694444445, 0.596, 0.5659999999999999, 0.131, 0.703, 0.594, 0, 0.673, 2.927, 1.311, 0.574, 0.838, 0.092, 0.061, 0.781, 0.09, 0.143, 2.22, 0.626, 2.165, 0.66, 0.638, 1.735, 2.278, 0.128, 0.419, 0.467, 0.214, 0.476, 0.5600000000000001, 0.8090000000000001, 2.624, 0.591, 0.381)
    3 = @(45086.51388888889, 20.495, 15.411, 0.729, 44.222, 36.746, 15.649, 53.774, 25.504, 11.844, 16.752, 18.027, 18.716, 5.149, 16.317, 22.575, 12.968, 1.319, 0.73, 238.235, 44.446, 15.053, 30.814, 16.608, 2.124, 27.869, 13.163, 11.444, 13.592, 18.89, 0.343, 49.44, 8.744999999999999, 18.38)
    4 = @(45086.52083333334, 22.935, 17.239, 0.792, 49.688, 41.226, 17.723, 69.82899999999999, 28.205, 12.984, 18.727, 20.126, 21.063, 5.779, 18.175, 25.459, 14.622, 0.981, 0.705, 267.092, 50.094, 16.779, 34.41, 18.259, 2.368, 34.385, 14.707, 12.843, 15.207, 21.135, 0.223, 64.095, 9.679, 20.63)
    5 = @(45086.52777777778, 10.95, 8.23, 0.38, 23.68, 19.69, 8.369999999999999, 37.53, 13.58, 6.36, 8.970000000000001, 9.640000000000001, 10.03, 2.77, 8.76, 12.12, 6.9, 0.7, 0.32, 124.11, 24.02, 8.09, 16.63, 8.83, 1.13, 17.87, 7.05, 6.15, 7.31, 10.13, 0.17, 34.64, 4.69, 9.83)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, $i + 1).Value = $vals[$i]
    }
}

# --- 3. Drop the now-unused trailing row 6 (also true in xlsx: rows used to
#     run to 6, new data only needs 5). ---
$ws.Rows.Item(6).Delete()

